$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new doctor record (D003, Xu Shuwei, Doctor, M, 21) is inserted at
# row 2, pushing the previous row-2/row-3 records down to rows 3/4.
# Row 2 (new): D003, Xu Shuwei, Doctor, M, 21
# Row 3 (was row 2 data's neighbour D002): D002, Emily Clarke, Doctor, Female, 38.0
# Row 4 (was row 2): D001, John Smith, Doctor, Male, 45.0

$ws.Cells.Item(2,1).Value = "D003"
$ws.Cells.Item(2,2).Value = "Xu Shuwei"
$ws.Cells.Item(2,3).Value = "Doctor"
$ws.Cells.Item(2,4).Value = "M"

$ws.Cells.Item(3,1).Value = "D002"
$ws.Cells.Item(3,2).Value = "Emily Clarke"
$ws.Cells.Item(3,3).Value = "Doctor"
$ws.Cells.Item(3,4).Value = "Female"

$ws.Cells.Item(4,1).Value = "D001"
$ws.Cells.Item(4,2).Value = "John Smith"
$ws.Cells.Item(4,3).Value = "Doctor"
$ws.Cells.Item(4,4).Value = "Male"

# The Age column values ("21", "38.0", "45.0") need to be stored as TEXT
# (shared-string) cells rather than numbers, matching the source data.
# Writing a numeric-looking string straight into Value would turn it into
# a genuine number, and forcing text via a leading apostrophe or a "@"
# number format on the destination cell leaves a stray quote-prefix /
# text style behind. Building the text value in an out-of-the-way scratch
# cell and pasting *values only* into the destination keeps the
# destination cell's style untouched (default style).
$scratch = $ws.Cells.Item(100,100)
$scratch.NumberFormat = "@"

$scratch.Value = "21"
$scratch.Copy()
$ws.Cells.Item(2,5).PasteSpecial(-4163)

$scratch.Value = "38.0"
$scratch.Copy()
$ws.Cells.Item(3,5).PasteSpecial(-4163)

$scratch.Value = "45.0"
$scratch.Copy()
$ws.Cells.Item(4,5).PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
